$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bestiary")
$lo = $ws.ListObjects.Item("Table2")

Write-Output "ROW HEIGHTS BEFORE:"
for ($r=1; $r -le 20; $r++) {
  $h = $ws.Rows.Item($r).RowHeight
  $v = $ws.Cells.Item($r,1).Value()
  Write-Output ("$r : $v h=$h")
}

$ws.Cells.Item(20,1).Value = "Giant Weasel"
$ws.Cells.Item(20,2).Value = "A giant weasel"
$ws.Cells.Item(20,3).Value = 11
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,6).Value = "Attack: 1d6 damage"

$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("A2:A28"), 0, 1)
$lo.Sort.Header = 0
$lo.Sort.Apply()

Write-Output "ROW HEIGHTS AFTER SORT:"
for ($r=1; $r -le 22; $r++) {
  $h = $ws.Rows.Item($r).RowHeight
  $v = $ws.Cells.Item($r,1).Value()
  Write-Output ("$r : $v h=$h")
}
